$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.401.84'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.847.65'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D5").Value = '''240.30'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '''0.6290'
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D8").Value = '''0.07610'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").Value = '''0.2932'
$ws.Range("D10").Value = '''24.48'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '''0.07744'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '1.851.23'
$ws.Range("E12").Value = '  -6.74%  '
$ws.Range("D13").Value = '''5.003'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = '''0.00001090'
$ws.Range("E14").Value = '  +9.60%  '
$ws.Range("D15").Value = '''0.6797'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '''83.77'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.098.84'
$ws.Range("E17").Value = '  -7.31%  '
$ws.Range("D18").Value = '''6.198'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = '29.417.10'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '''228.87'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D23").Value = '''7.468'
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D25").Value = '''157.40'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").Value = '''0.1398'
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("D27").Value = '''8.362'
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = '''1.301'
$ws.Range("E30").Value = '  +4.27%  '
$ws.Range("D31").Value = '''0.05584'
$ws.Range("E31").Value = '  -2.20%  '
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").Value = '''1.848'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = '''1.156'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = '''0.7104'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '1.233.68'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01800'
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.772'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").Value = '''6.438'
$ws.Range("E41").Value = '  +5.66%  '
$ws.Range("D42").Value = '''0.9077'
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").Value = '''101.90'
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").Value = '''66.13'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = '''0.00000000120'
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("D47").Value = '''7.174'
$ws.Range("D48").Value = '''0.4025'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").Value = '''8.986'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").Value = '''1.679'
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("E51").Value = '  -0.52%  '

# Reset style on cells that needed a quote-prefix so no stray style/quotePrefix
# flag is left behind on cells that must remain plain "General" text.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
